# Add a source-reference entry (hyperlinked URL + retrieval date) to the
# "Quellen" content placeholder on the slide identified by SlideID 279.
#
# Target paragraph is the first otherwise-empty trailing paragraph of the
# numbered source list; two runs are inserted before its (empty) end:
#   1) the hyperlinked URL
#   2) a plain-text " <date>, <time> Uhr" suffix

$p = $ppt.ActivePresentation

# Locate the slide by its persistent SlideID (robust to reordering).
$slide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq 279) {
        $slide = $candidate
    }
}

# Locate the content placeholder shape by its shape Id.
$shape = $null
for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
    $candidateShape = $slide.Shapes.Item($j)
    if ($candidateShape.Id -eq 2) {
        $shape = $candidateShape
    }
}

$textRange = $shape.TextFrame.TextRange

# Find the first empty paragraph (the blank numbered line meant to hold the
# new source entry).
$targetIndex = 0
for ($k = 1; $k -le $textRange.Paragraphs().Count; $k++) {
    if ($targetIndex -eq 0 -and $textRange.Paragraphs($k).Text.Trim() -eq "") {
        $targetIndex = $k
    }
}

$url = "https://msdn.microsoft.com/de-de/library/dn133186(v=sql.120).aspx"
$suffix = " 15.01.2018, 17.03 Uhr"

$paragraph = $textRange.Paragraphs($targetIndex)

# Insert the combined text as a single run first, then carve the hyperlink
# out of just the URL portion so the date/time suffix stays a plain run.
$combined = $paragraph.InsertBefore($url + $suffix)
$combined.Font.Size = 14

$paragraph = $textRange.Paragraphs($targetIndex)
$urlRange = $textRange.Characters($paragraph.Start, $url.Length)
$urlRange.ActionSettings(1).Hyperlink.Address = $url
